$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.5871
$ws.Range("D3").Value = 0.3037
$ws.Range("E3").Value = 0.0566
$ws.Range("F3").Value = 0.0023
$ws.Range("G3").Value = 0.0315
$ws.Range("H3").Value = 0.0187

[void]$ws.Range("J3").Select()
